# Update with Correct Forecast output
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Forecast Comparison
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column B ("Week_Start_Date") shifting ASIN / MyForecast /
# Amazon forecasts / Product Title / is_holiday_week one column to the right.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Final data (Week label, Week start date, MyForecast, Amazon Mean,
# Amazon P70, Amazon P80, Amazon P90) for rows 2-17.
$data = @(
    @("W1",  "2025-01-05", 198, 300, 357, 414, 500),
    @("W2",  "2025-01-12", 188, 259, 314, 375, 471),
    @("W3",  "2025-01-19", 197, 254, 309, 373, 476),
    @("W4",  "2025-01-26", 202, 242, 292, 346, 433),
    @("W5",  "2025-02-02", 178, 168, 203, 241, 301),
    @("W6",  "2025-02-09", 169, 165, 199, 236, 295),
    @("W7",  "2025-02-16", 196, 161, 195, 234, 296),
    @("W8",  "2025-02-23", 193, 167, 203, 246, 315),
    @("W9",  "2025-03-02", 199, 162, 195, 232, 289),
    @("W10", "2025-03-09", 196, 166, 201, 242, 306),
    @("W11", "2025-03-16", 191, 156, 191, 237, 311),
    @("W12", "2025-03-23", 184, 156, 190, 235, 307),
    @("W13", "2025-03-30", 184, 155, 188, 227, 290),
    @("W14", "2025-04-06", 180, 146, 178, 217, 280),
    @("W15", "2025-04-13", 181, 141, 172, 210, 271),
    @("W16", "2025-04-20", 176, 141, 171, 209, 269)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]

    # Column A: Week label (e.g. "W01" -> "W1")
    $ws.Cells.Item($r, 1).Value = $row[0]

    # Column B: Week_Start_Date, stored as literal text (not an Excel date).
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]

    # Column D: MyForecast (refreshed values)
    $ws.Cells.Item($r, 4).Value = $row[2]

    # Columns E-H: Amazon Mean / P70 / P80 / P90 forecasts
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]

    # Column J: is_holiday_week now stored as a boolean instead of a number.
    $ws.Cells.Item($r, 10).Value = $false
}

# ---------------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------------
# These numeric-looking values must stay stored as text, matching the rest
# of the column, so prefix with a literal quote to force text entry.
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(9, 2).Value  = "'3012"
$summary.Cells.Item(10, 2).Value = "'1520"
$summary.Cells.Item(12, 2).Value = "'202"
$summary.Cells.Item(14, 2).Value = "'169"
